# Update predidx (D) and pred_name (E) columns for specific rows
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 3;  D = "[1, 0, 0, 0, 0, 0, 0]"; E = "['Normal']" },
    @{ Row = 5;  D = "[1, 0, 0, 0, 0, 0, 0]"; E = "['Normal']" },
    @{ Row = 6;  D = "[1, 0, 0, 0, 0, 0, 0]"; E = "['Normal']" },
    @{ Row = 11; D = "[1, 0, 1, 0, 1, 0, 0]"; E = "['Normal', 'HardwareFault', 'RegulationViolation']" },
    @{ Row = 12; D = "[1, 0, 1, 0, 0, 0, 0]"; E = "['Normal', 'HardwareFault']" },
    @{ Row = 15; D = "[1, 0, 0, 0, 0, 0, 0]"; E = "['Normal']" },
    @{ Row = 25; D = "[1, 0, 0, 0, 0, 0, 0]"; E = "['Normal']" },
    @{ Row = 26; D = "[0, 0, 0, 0, 0, 0, 0]"; E = "[]" },
    @{ Row = 29; D = "[1, 0, 0, 0, 0, 0, 1]"; E = "['Normal', 'SoftwareFault']" },
    @{ Row = 36; D = "[1, 1, 1, 0, 0, 0, 0]"; E = "['Normal', 'SurroundingEnvironment', 'HardwareFault']" },
    @{ Row = 38; D = "[1, 0, 0, 0, 0, 0, 1]"; E = "['Normal', 'SoftwareFault']" },
    @{ Row = 54; D = "[1, 0, 0, 0, 0, 0, 0]"; E = "['Normal']" },
    @{ Row = 61; D = "[1, 0, 0, 0, 0, 0, 1]"; E = "['Normal', 'SoftwareFault']" },
    @{ Row = 71; D = "[1, 0, 0, 0, 0, 0, 0]"; E = "['Normal']" },
    @{ Row = 81; D = "[1, 0, 0, 0, 0, 0, 0]"; E = "['Normal']" },
    @{ Row = 84; D = "[1, 0, 0, 0, 0, 0, 0]"; E = "['Normal']" },
    @{ Row = 89; D = "[1, 0, 0, 0, 0, 0, 0]"; E = "['Normal']" },
    @{ Row = 97; D = "[1, 0, 0, 0, 0, 0, 0]"; E = "['Normal']" }
)

foreach ($u in $updates) {
    $ws.Range("D" + $u.Row).Value = $u.D
    $ws.Range("E" + $u.Row).Value = $u.E
}
